# "version for hazard workshop"
# Refresh the SPO2IDA input parameters (column B, rows 2-8) with the values
# used for the hazard-workshop run. The dependent formulas in B10:B20
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.104435862
$ws.Range("B3").Value = 0.233993594
$ws.Range("B4").Value = 0.89770978700000004
$ws.Range("B5").Value = 0.94027882799999996
$ws.Range("B6").Value = 2093.9
$ws.Range("B7").Value = 2093.9
$ws.Range("B8").Value = 91.434464939999998

# Restore the view the author left the sheet in: scrolled down a bit with
# B19 (the last "r" / ac row of the GET PARAMETERS block) selected.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B19").Select()
